# Adds Test Data for Hungary/Russia/Finland Market.
#
# Strategy: the three new "Market" worksheets (Russia, Finland, Hungary) use
# exactly the same cell styling / layout as the existing Spain / Croatia /
# Greece sheets (same cellXfs indices, same merged header cell C1:D1, same
# "accessory list" rows) except that the legacy duplicated "RDS800" row
# (row 16, a second copy of row 7) is not present - the new sheets only
# have 17 rows instead of 18.
#
# So for each new market: copy the "Spain" sheet (it already carries the
# right style indices), rename it, overwrite the market name / NGC code,
# drop the duplicated RDS800 row, and select the sheet's used range (mirrors
# how the sheets look when freshly created from a template and then tidied
# up). The last sheet created (Hungary) is left as the active / selected
# tab, matching the workbook's final `activeTab`.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Spain")

function Add-MarketSheet($wb, $template, $sheetName, $ngcCode, $marketName) {
    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $afterSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $sheetName

    # Fill in the NGC code (B3) before the market name (B1) - matches the
    # order the new shared strings were appended to the workbook.
    $newSheet.Range("B3").Value = $ngcCode
    $newSheet.Range("B1").Value = $marketName

    # Remove the legacy duplicate "RDS800" row that Spain/Croatia/Greece
    # carry at row 16 - the new sheets don't have it (18 rows -> 17 rows).
    $newSheet.Rows.Item(16).Delete()

    return $newSheet
}

$russia = Add-MarketSheet $wb $template "Russia" "NGC-2929/T2925" "Russia Market"
[void]$russia.Range("A1:D17").Select()

$finland = Add-MarketSheet $wb $template "Finland" "NGC-3130/T2957" "Finland Market"
[void]$finland.Range("A1:D17").Select()

$hungary = Add-MarketSheet $wb $template "Hungary" "NGC-3104/T3006" "Hungary Market"
[void]$hungary.Range("L11").Select()
[void]$hungary.Select()
